$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.356.48'
$ws.Range("E2").Value = '  +5.20%  '

# Row 3
$ws.Range("D3").Value = '1.816.33'
$ws.Range("E3").Value = '  +5.37%  '

# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = '''318.28'
$ws.Range("E5").Value = '  +2.66%  '

# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.14%  '

# Row 7
$ws.Range("D7").Value = '''0.5748'
$ws.Range("E7").Value = '  +18.85%  '

# Row 8
$ws.Range("D8").Value = '''0.3866'
$ws.Range("E8").Value = '  +11.47%  '

# Row 9
$ws.Range("D9").Value = '''0.07622'
$ws.Range("E9").Value = '  +5.40%  '

# Row 10
$ws.Range("D10").Value = '''43.19'
$ws.Range("E10").Value = '  -0.01%  '

# Row 11
$ws.Range("E11").Value = '  +8.69%  '

# Row 12
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  -0.08%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '''21.23'
$ws.Range("E13").Value = '  +7.17%  '

# Row 14
$ws.Range("D14").Value = '''6.251'
$ws.Range("E14").Value = '  +6.64%  '

# Row 15
$ws.Range("D15").Value = '1.810.76'
$ws.Range("E15").Value = '  +4.90%  '

# Row 16
$ws.Range("D16").Value = '''7.310'

# Row 17
$ws.Range("D17").Value = '''92.18'
$ws.Range("E17").Value = '  +6.17%  '

# Row 18
$ws.Range("D18").Value = '''0.00001076'
$ws.Range("E18").Value = '  +4.15%  '

# Row 19
$ws.Range("D19").Value = '''0.06469'
$ws.Range("E19").Value = '  +1.11%  '

# Row 20
$ws.Range("D20").Value = '''0.9998'
$ws.Range("E20").Value = '  -0.12%  '

# Row 21
$ws.Range("E21").Value = '  +4.79%  '

# Row 22
$ws.Range("D22").Value = '''5.990'
$ws.Range("E22").Value = '  +5.23%  '

# Row 23
$ws.Range("D23").Value = '28.377.28'
$ws.Range("E23").Value = '  +5.03%  '

# Row 24
$ws.Range("D24").Value = '''11.31'
$ws.Range("E24").Value = '  +3.54%  '

# Row 25
$ws.Range("E25").Value = '  +1.54%  '

# Row 26
$ws.Range("D26").Value = '''20.84'
$ws.Range("E26").Value = '  +4.90%  '

# Row 27
$ws.Range("D27").Value = '''157.68'
$ws.Range("E27").Value = '  +2.50%  '

# Row 28
$ws.Range("D28").Value = '''2.408'
$ws.Range("E28").Value = '  +16.72%  '

# Row 29
$ws.Range("D29").Value = '2.023.88'
$ws.Range("E29").Value = '  +5.01%  '

# Row 30
$ws.Range("D30").Value = '''123.65'
$ws.Range("E30").Value = '  +2.39%  '

# Row 31
$ws.Range("D31").Value = '''1.175'
$ws.Range("E31").Value = '  +14.12%  '

# Row 32
$ws.Range("D32").Value = '''0.1062'
$ws.Range("E32").Value = '  +14.22%  '

# Row 33
$ws.Range("D33").Value = '''5.773'
$ws.Range("E33").Value = '  +7.55%  '

# Row 34
$ws.Range("D34").Value = '''3.631'
$ws.Range("E34").Value = '  +0.21%  '

# Row 35
$ws.Range("D35").Value = '''8.879'
$ws.Range("E35").Value = '  +19.19%  '

# Row 36
$ws.Range("D36").Value = '''0.02315'
$ws.Range("E36").Value = '  +6.73%  '

# Row 37
$ws.Range("D37").Value = '''0.2165'
$ws.Range("E37").Value = '  +9.00%  '

# Row 38
$ws.Range("D38").Value = '''11.68'
$ws.Range("E38").Value = '  +7.21%  '

# Row 39
$ws.Range("E39").Value = '  +7.62%  '

# Row 40
$ws.Range("D40").Value = '''0.06085'
$ws.Range("E40").Value = '  +2.47%  '

# Row 41
$ws.Range("D41").Value = '''5.031'
$ws.Range("E41").Value = '  +6.59%  '

# Row 42
$ws.Range("D42").Value = '''1.000'
$ws.Range("E42").Value = '  -0.06%  '

# Row 43
$ws.Range("E43").Value = '  +3.61%  '

# Row 44
$ws.Range("D44").Value = '''1.379'

# Row 45
$ws.Range("D45").Value = '''13.39'
$ws.Range("E45").Value = '  +4.50%  '

# Row 46
$ws.Range("D46").Value = '''0.5979'
$ws.Range("E46").Value = '  +6.98%  '

# Row 47
$ws.Range("D47").Value = '''3.704'
$ws.Range("E47").Value = '  +3.58%  '

# Row 48
$ws.Range("D48").Value = '''122.07'
$ws.Range("E48").Value = '  +2.59%  '

# Row 49
$ws.Range("D49").Value = '''1.941'
$ws.Range("E49").Value = '  +5.30%  '

# Row 50
$ws.Range("D50").Value = '''1.147'
$ws.Range("E50").Value = '  +4.44%  '

# Row 51
$ws.Range("D51").Value = '''0.06841'
$ws.Range("E51").Value = '  +3.04%  '
